# Auto commit at 2025-12-13 17:24:51.28
# Refresh the "Metrics" sheet's source values (B2:B13) and let the
# dependent formulas on the "today" sheet (B11:B22, E11:E22, F11:F22,
# which reference Metrics!B2:B13) recalculate naturally.

$wb = $excel.ActiveWorkbook
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 159471.35999999999
$metrics.Range("B3").Value  = 136928.12000000002
$metrics.Range("B4").Value  = 49002.239999999998
$metrics.Range("B5").Value  = 6560
$metrics.Range("B6").Value  = 5362178.4700000007
$metrics.Range("B7").Value  = 4537281.080000001
$metrics.Range("B8").Value  = 1580959.1200000003
$metrics.Range("B9").Value  = 209267
$metrics.Range("B10").Value = 33827559.459999993
$metrics.Range("B11").Value = 31812556.240000002
$metrics.Range("B12").Value = 11862681.159999995
$metrics.Range("B13").Value = 1306897

# Move the cell selection on the Metrics sheet to D20. Activate the sheet
# first so the Range.Select() call is valid, then restore the workbook's
# originally active sheet ("today") so the overall active-tab/view state
# is left untouched.
$originalActive = $wb.ActiveSheet
$metrics.Activate()
$metrics.Range("D20").Select()
$originalActive.Activate()
